$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values that changed, keeping them stored as text
$ws.Range("D2").Value = "'242.53"
$ws.Range("D3").Value = "'22.16"
$ws.Range("D4").Value = "'5.425"
$ws.Range("D5").Value = "'0.05748"
$ws.Range("D6").Value = "'3.431"
$ws.Range("D8").Value = "'0.8116"
$ws.Range("D9").Value = "'0.8612"
$ws.Range("D10").Value = "'0.1443"
$ws.Range("D11").Value = "'0.07323"
$ws.Range("D12").Value = "'0.03063"
$ws.Range("D13").Value = "'0.03112"
$ws.Range("D14").Value = "'0.09396"
$ws.Range("D15").Value = "'3.935"
$ws.Range("D16").Value = "'0.001593"
$ws.Range("D17").Value = "'0.04846"
$ws.Range("D18").Value = "'0.0005856"
$ws.Range("D19").Value = "'0.006412"
$ws.Range("D21").Value = "'0.0009973"
$ws.Range("D23").Value = "'3.722"
$ws.Range("D24").Value = "'2.191"
$ws.Range("D25").Value = "'0.3269"
$ws.Range("D27").Value = "'0.0004004"
$ws.Range("D40").Value = "'0.03859"
$ws.Range("D41").Value = "'0.006740"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("D43").Value = "'0.002423"
$ws.Range("D44").Value = "'0.006452"
$ws.Range("D45").Value = "'0.00005598"
$ws.Range("D47").Value = "'0.3804"
$ws.Range("D48").Value = "'0.1448"
$ws.Range("D49").Value = "'0.00002102"

# Update Hora column (G) from 13 to 14 for all data rows, keeping text type
$ws.Range("G2").Value = "'14"
$ws.Range("G3").Value = "'14"
$ws.Range("G4").Value = "'14"
$ws.Range("G5").Value = "'14"
$ws.Range("G6").Value = "'14"
$ws.Range("G7").Value = "'14"
$ws.Range("G8").Value = "'14"
$ws.Range("G9").Value = "'14"
$ws.Range("G10").Value = "'14"
$ws.Range("G11").Value = "'14"
$ws.Range("G12").Value = "'14"
$ws.Range("G13").Value = "'14"
$ws.Range("G14").Value = "'14"
$ws.Range("G15").Value = "'14"
$ws.Range("G16").Value = "'14"
$ws.Range("G17").Value = "'14"
$ws.Range("G18").Value = "'14"
$ws.Range("G19").Value = "'14"
$ws.Range("G20").Value = "'14"
$ws.Range("G21").Value = "'14"
$ws.Range("G22").Value = "'14"
$ws.Range("G23").Value = "'14"
$ws.Range("G24").Value = "'14"
$ws.Range("G25").Value = "'14"
$ws.Range("G26").Value = "'14"
$ws.Range("G27").Value = "'14"
$ws.Range("G28").Value = "'14"
$ws.Range("G29").Value = "'14"
$ws.Range("G30").Value = "'14"
$ws.Range("G31").Value = "'14"
$ws.Range("G32").Value = "'14"
$ws.Range("G33").Value = "'14"
$ws.Range("G34").Value = "'14"
$ws.Range("G35").Value = "'14"
$ws.Range("G36").Value = "'14"
$ws.Range("G37").Value = "'14"
$ws.Range("G38").Value = "'14"
$ws.Range("G39").Value = "'14"
$ws.Range("G40").Value = "'14"
$ws.Range("G41").Value = "'14"
$ws.Range("G42").Value = "'14"
$ws.Range("G43").Value = "'14"
$ws.Range("G44").Value = "'14"
$ws.Range("G45").Value = "'14"
$ws.Range("G46").Value = "'14"
$ws.Range("G47").Value = "'14"
$ws.Range("G48").Value = "'14"
$ws.Range("G49").Value = "'14"
$ws.Range("G50").Value = "'14"
$ws.Range("G51").Value = "'14"
